$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact literal text value (avoids Excel
# auto-converting numeric-looking strings like "1.001" into numbers,
# and resets the style back to Normal so no quote-prefix style sticks).
function Set-ExactText($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-ExactText $ws.Range("D2") '28.407.13'
Set-ExactText $ws.Range("E2") '  +4.44%  '
Set-ExactText $ws.Range("D3") '1.800.50'
Set-ExactText $ws.Range("E3") '  +1.70%  '
Set-ExactText $ws.Range("E4") '  -0.22%  '
Set-ExactText $ws.Range("D5") '314.89'
Set-ExactText $ws.Range("E5") '  +0.47%  '
Set-ExactText $ws.Range("D6") '1.001'
Set-ExactText $ws.Range("E6") '  -0.07%  '
Set-ExactText $ws.Range("D7") '0.5505'
Set-ExactText $ws.Range("E7") '  +4.64%  '
Set-ExactText $ws.Range("D8") '0.3851'
Set-ExactText $ws.Range("E8") '  +5.32%  '
Set-ExactText $ws.Range("D9") '0.07598'
Set-ExactText $ws.Range("E9") '  +3.49%  '
Set-ExactText $ws.Range("D10") '42.54'
Set-ExactText $ws.Range("E10") '  -0.66%  '
Set-ExactText $ws.Range("D11") '1.128'
Set-ExactText $ws.Range("E11") '  +3.88%  '
Set-ExactText $ws.Range("E12") '  -0.26%  '
Set-ExactText $ws.Range("D13") '21.17'
Set-ExactText $ws.Range("E13") '  +4.15%  '
Set-ExactText $ws.Range("D14") '6.184'
Set-ExactText $ws.Range("E14") '  +2.30%  '
Set-ExactText $ws.Range("D15") '7.412'
Set-ExactText $ws.Range("E15") '  +7.06%  '
Set-ExactText $ws.Range("D16") '1.804.73'
Set-ExactText $ws.Range("E16") '  +2.06%  '
Set-ExactText $ws.Range("D17") '91.95'
Set-ExactText $ws.Range("E17") '  +3.79%  '
Set-ExactText $ws.Range("D18") '0.00001073'
Set-ExactText $ws.Range("E18") '  +3.00%  '
Set-ExactText $ws.Range("D19") '0.06448'
Set-ExactText $ws.Range("E19") '  +0.46%  '
Set-ExactText $ws.Range("D20") '1.000'
Set-ExactText $ws.Range("E20") '  -0.16%  '
Set-ExactText $ws.Range("D21") '17.37'
Set-ExactText $ws.Range("E21") '  +4.25%  '
Set-ExactText $ws.Range("D22") '5.972'
Set-ExactText $ws.Range("E22") '  +2.81%  '
Set-ExactText $ws.Range("D23") '28.414.14'
Set-ExactText $ws.Range("E23") '  +4.03%  '
Set-ExactText $ws.Range("D24") '11.44'
Set-ExactText $ws.Range("E24") '  +1.77%  '
Set-ExactText $ws.Range("D25") '2.127'
Set-ExactText $ws.Range("E25") '  +1.02%  '
Set-ExactText $ws.Range("D26") '158.76'
Set-ExactText $ws.Range("E26") '  +2.96%  '
Set-ExactText $ws.Range("D27") '20.71'
Set-ExactText $ws.Range("E27") '  +3.28%  '
Set-ExactText $ws.Range("D28") '2.406'
Set-ExactText $ws.Range("E28") '  +3.44%  '
Set-ExactText $ws.Range("D29") '2.013.18'
Set-ExactText $ws.Range("E29") '  +2.12%  '
Set-ExactText $ws.Range("D30") '123.79'
Set-ExactText $ws.Range("E30") '  +2.58%  '
Set-ExactText $ws.Range("D31") '1.120'
Set-ExactText $ws.Range("E31") '  +5.97%  '
Set-ExactText $ws.Range("D32") '0.1026'
Set-ExactText $ws.Range("E32") '  +5.64%  '
Set-ExactText $ws.Range("D33") '5.772'
Set-ExactText $ws.Range("E33") '  +4.17%  '
Set-ExactText $ws.Range("D34") '3.693'
Set-ExactText $ws.Range("E34") '  +1.96%  '
Set-ExactText $ws.Range("D35") '0.2309'
Set-ExactText $ws.Range("E35") '  +14.66%  '
Set-ExactText $ws.Range("D36") '0.06400'
Set-ExactText $ws.Range("E36") '  +7.70%  '
Set-ExactText $ws.Range("E37") '  +4.95%  '
Set-ExactText $ws.Range("D38") '5.150'
Set-ExactText $ws.Range("E38") '  +6.85%  '
Set-ExactText $ws.Range("D39") '8.777'
Set-ExactText $ws.Range("E39") '  +9.20%  '
Set-ExactText $ws.Range("E40") '  +4.35%  '
Set-ExactText $ws.Range("E41") '  +5.12%  '
Set-ExactText $ws.Range("D42") '1.000'
Set-ExactText $ws.Range("E42") '  -0.02%  '
Set-ExactText $ws.Range("D43") '1.158'
Set-ExactText $ws.Range("E43") '  +2.16%  '
Set-ExactText $ws.Range("E44") '  -3.42%  '
Set-ExactText $ws.Range("D45") '13.62'
Set-ExactText $ws.Range("E45") '  +3.87%  '
Set-ExactText $ws.Range("D46") '0.5975'
Set-ExactText $ws.Range("E46") '  +4.40%  '
Set-ExactText $ws.Range("D47") '3.679'
Set-ExactText $ws.Range("E47") '  +1.58%  '
Set-ExactText $ws.Range("D48") '126.53'
Set-ExactText $ws.Range("E48") '  +4.91%  '
Set-ExactText $ws.Range("D49") '1.981'
Set-ExactText $ws.Range("E49") '  +5.39%  '
Set-ExactText $ws.Range("D50") '1.151'
Set-ExactText $ws.Range("E50") '  +3.93%  '
Set-ExactText $ws.Range("D51") '0.06905'
Set-ExactText $ws.Range("E51") '  +2.84%  '
